# Auto-generated update of market price / profit figures across multiple sheets
# (currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 327.33334
$ws.Range("I28").Value = 480
$ws.Range("J28").Value = 136.5
$ws.Range("K28").Value = 480
$ws.Range("L28").Value = 136.5
$ws.Range("M28").Value = 5
$ws.Range("N28").Value = -1106.5

$ws.Range("H86").Value = 23137.2
$ws.Range("I86").Value = 3733.3333
$ws.Range("J86").Value = 52243
$ws.Range("K86").Value = 3733.3333
$ws.Range("L86").Value = 52243
$ws.Range("M86").Value = -2610.3333
$ws.Range("N86").Value = -54489

$ws.Range("H89").Value = 23137.2
$ws.Range("I89").Value = 3733.3333
$ws.Range("J89").Value = 52243
$ws.Range("K89").Value = 18666.6665
$ws.Range("L89").Value = 261215
$ws.Range("M89").Value = -13050.6665
$ws.Range("N89").Value = -272447

$ws.Range("H98").Value = 1550
$ws.Range("I98").Value = 825
$ws.Range("J98").Value = 3000
$ws.Range("K98").Value = 825
$ws.Range("L98").Value = 3000
$ws.Range("M98").Value = 673
$ws.Range("N98").Value = -5996

$ws.Range("H106").Value = 1685.8334
$ws.Range("I106").Value = 1302.7778
$ws.Range("K106").Value = 1302.7778
$ws.Range("M106").Value = -671.7778000000001

$ws.Range("H122").Value = 1550
$ws.Range("I122").Value = 825
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 2475
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -25
$ws.Range("N122").Value = -13900

$ws.Range("H129").Value = 861.5325
$ws.Range("J129").Value = 874.61646
$ws.Range("L129").Value = 2623.84938
$ws.Range("N129").Value = -12623.84938

$ws.Range("H138").Value = 2503.6792
$ws.Range("J138").Value = 3361.8108
$ws.Range("L138").Value = 10085.4324
$ws.Range("N138").Value = -20365.4324

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2175.3572
$ws.Range("I63").Value = 2171.25
$ws.Range("J63").Value = 2200
$ws.Range("K63").Value = 2171.25
$ws.Range("L63").Value = 2200
$ws.Range("M63").Value = -1485.25
$ws.Range("N63").Value = -3572

$ws.Range("H66").Value = 2175.3572
$ws.Range("I66").Value = 2171.25
$ws.Range("J66").Value = 2200
$ws.Range("K66").Value = 10856.25
$ws.Range("L66").Value = 11000
$ws.Range("M66").Value = -7424.25
$ws.Range("N66").Value = -17864

$ws.Range("H88").Value = 50431.145
$ws.Range("I88").Value = 1581.2
$ws.Range("K88").Value = 1581.2
$ws.Range("M88").Value = -1175.2

$ws.Range("H91").Value = 50431.145
$ws.Range("I91").Value = 1581.2
$ws.Range("K91").Value = 1581.2
$ws.Range("M91").Value = -177.2

$ws.Range("H132").Value = 15951.333
$ws.Range("I132").Value = 1547.55
$ws.Range("J132").Value = 33956.062
$ws.Range("K132").Value = 4642.65
$ws.Range("L132").Value = 101868.186
$ws.Range("M132").Value = -2112.65
$ws.Range("N132").Value = -106928.186

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 31185.646
$ws.Range("J86").Value = 2792.6667
$ws.Range("L86").Value = 2792.6667
$ws.Range("N86").Value = -5038.6667

$ws.Range("H89").Value = 31185.646
$ws.Range("J89").Value = 2792.6667
$ws.Range("L89").Value = 13963.3335
$ws.Range("N89").Value = -25195.3335

$ws.Range("H94").Value = 3672
$ws.Range("I94").Value = 826
$ws.Range("J94").Value = 6873.75
$ws.Range("K94").Value = 826
$ws.Range("L94").Value = 6873.75
$ws.Range("M94").Value = -375
$ws.Range("N94").Value = -7775.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 1000
$ws.Range("J4").Value = 1000
$ws.Range("L4").Value = 1000
$ws.Range("N4").Value = -1224

$ws.Range("H62").Value = 6002.25
$ws.Range("I62").Value = 5000
$ws.Range("J62").Value = 6336.3335
$ws.Range("K62").Value = 5000
$ws.Range("L62").Value = 6336.3335
$ws.Range("M62").Value = -4376
$ws.Range("N62").Value = -7584.3335

$ws.Range("H65").Value = 6002.25
$ws.Range("I65").Value = 5000
$ws.Range("J65").Value = 6336.3335
$ws.Range("K65").Value = 25000
$ws.Range("L65").Value = 31681.6675
$ws.Range("M65").Value = -21880
$ws.Range("N65").Value = -37921.6675

$ws.Range("H122").Value = 881.73334
$ws.Range("I122").Value = 940.7273
$ws.Range("J122").Value = 719.5
$ws.Range("K122").Value = 2822.1819
$ws.Range("L122").Value = 2158.5
$ws.Range("M122").Value = -372.1819
$ws.Range("N122").Value = -7058.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 639.5641000000001
$ws.Range("I5").Value = 464.5
$ws.Range("K5").Value = 1393.5
$ws.Range("M5").Value = -1281.5

$ws.Range("H68").Value = 1293.1818
$ws.Range("J68").Value = 1438.8334
$ws.Range("L68").Value = 4316.5002
$ws.Range("N68").Value = -5938.5002

$ws.Range("H71").Value = 1293.1818
$ws.Range("J71").Value = 1438.8334
$ws.Range("L71").Value = 12949.5006
$ws.Range("N71").Value = -21061.5006

$ws.Range("H107").Value = 4827
$ws.Range("I107").Value = 9580.546
$ws.Range("J107").Value = 804.7692
$ws.Range("K107").Value = 28741.638
$ws.Range("L107").Value = 2414.3076
$ws.Range("M107").Value = -26821.638
$ws.Range("N107").Value = -6254.3076

$ws.Range("H131").Value = 104999.9
$ws.Range("J131").Value = 110715.34
$ws.Range("L131").Value = 332146.02
$ws.Range("N131").Value = -342226.02

$ws.Range("H135").Value = 639.5641000000001
$ws.Range("I135").Value = 464.5
$ws.Range("K135").Value = 4180.5
$ws.Range("M135").Value = -1645.5

$ws.Range("H140").Value = 5619.72
$ws.Range("I140").Value = 7319.933
$ws.Range("J140").Value = 3069.4
$ws.Range("K140").Value = 21959.799
$ws.Range("L140").Value = 9208.200000000001
$ws.Range("M140").Value = -16779.799
$ws.Range("N140").Value = -19568.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8899.166999999999
$ws.Range("I80").Value = 14825.625
$ws.Range("J80").Value = 4158
$ws.Range("K80").Value = 14825.625
$ws.Range("L80").Value = 4158
$ws.Range("M80").Value = -13827.625
$ws.Range("N80").Value = -6154

$ws.Range("H83").Value = 8899.166999999999
$ws.Range("I83").Value = 14825.625
$ws.Range("J83").Value = 4158
$ws.Range("K83").Value = 74128.125
$ws.Range("L83").Value = 20790
$ws.Range("M83").Value = -69136.125
$ws.Range("N83").Value = -30774

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 869.61536
$ws.Range("I46").Value = 750.4167
$ws.Range("K46").Value = 750.4167
$ws.Range("M46").Value = -562.4167

$ws.Range("H136").Value = 38528.645
$ws.Range("I136").Value = 56711.332
$ws.Range("K136").Value = 170133.996
$ws.Range("M136").Value = -167583.996

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5697.143
$ws.Range("I62").Value = 5440
$ws.Range("J62").Value = 5800
$ws.Range("K62").Value = 5440
$ws.Range("L62").Value = 5800
$ws.Range("M62").Value = -4816
$ws.Range("N62").Value = -7048

$ws.Range("H65").Value = 5697.143
$ws.Range("I65").Value = 5440
$ws.Range("J65").Value = 5800
$ws.Range("K65").Value = 27200
$ws.Range("L65").Value = 29000
$ws.Range("M65").Value = -24080
$ws.Range("N65").Value = -35240

$ws.Range("H81").Value = 1487.5714
$ws.Range("I81").Value = 1435.5
$ws.Range("J81").Value = 1800
$ws.Range("K81").Value = 2871
$ws.Range("L81").Value = 3600
$ws.Range("M81").Value = -1810
$ws.Range("N81").Value = -5722

$ws.Range("H84").Value = 1487.5714
$ws.Range("I84").Value = 1435.5
$ws.Range("J84").Value = 1800
$ws.Range("K84").Value = 14355
$ws.Range("L84").Value = 18000
$ws.Range("M84").Value = -9051
$ws.Range("N84").Value = -28608
